# Agile Gantt chart.xlsx - update the Gantt "Scrolling Increment" control
# (cell C7 on the "Color" sheet) so the chart scrolls forward to show
# October/November/December instead of March/April/May. All of the date
# headers (row 6), the date ruler (row 7) and the conditional Gantt-bar
# helper formulas recalculate automatically off of this single input, so
# only the one cell needs to be written - Excel's calc engine (via the
# IFERROR(Project_Start+Scrolling_Increment, TODAY()) chain) takes care
# of everything downstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Color")

# Scrolling_Increment named range -> Color!$C$7
$ws.Range("C7").Value = 20

# Refresh the view: zoom out and move the selection/scroll position.
$win = $excel.ActiveWindow
$win.Zoom = 55
$ws.Range("C8").Select()
